# "reclassificando os tweets, choro /2"
# Re-label a batch of tweets in the "Treinamento" sheet (column B holds the
# 0/1 classification), then leave the view scrolled/zoomed near the bottom
# of that sheet as the active tab (mirrors the author's last on-screen state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

# --- reclassified rows (B column: 0 <-> 1) ---
$ws.Range("B7").Value = 0
$ws.Range("B14").Value = 1
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B20").Value = 1
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B27").Value = 1
$ws.Range("B39").Value = 0
$ws.Range("B42").Value = 1
$ws.Range("B48").Value = 0
$ws.Range("B50").Value = 1
$ws.Range("B51").Value = 0
$ws.Range("B68").Value = 0
$ws.Range("B85").Value = 0
$ws.Range("B91").Value = 0
$ws.Range("B96").Value = 1
$ws.Range("B105").Value = 0
$ws.Range("B113").Value = 1
$ws.Range("B115").Value = 1
$ws.Range("B121").Value = 0
$ws.Range("B123").Value = 0
$ws.Range("B125").Value = 0
$ws.Range("B149").Value = 1
$ws.Range("B153").Value = 1
$ws.Range("B167").Value = 1
$ws.Range("B171").Value = 1
$ws.Range("B183").Value = 1
$ws.Range("B188").Value = 0
$ws.Range("B190").Value = 0
$ws.Range("B192").Value = 1
$ws.Range("B207").Value = 1
$ws.Range("B208").Value = 1
$ws.Range("B215").Value = 1
$ws.Range("B225").Value = 1
$ws.Range("B230").Value = 1
$ws.Range("B237").Value = 1
$ws.Range("B257").Value = 1
$ws.Range("B258").Value = 1
$ws.Range("B262").Value = 0
$ws.Range("B263").Value = 1
$ws.Range("B269").Value = 1
$ws.Range("B275").Value = 1
$ws.Range("B285").Value = 1
$ws.Range("B287").Value = 1
$ws.Range("B288").Value = 1
$ws.Range("B289").Value = 1
$ws.Range("B291").Value = 1

# --- final view state: "Treinamento" becomes the active/selected tab,
# zoomed to 115%, scrolled near the bottom, with A301 selected ---
$ws.Activate()
$ws.Select()
$excel.ActiveWindow.ScrollRow = 282
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 115
$ws.Range("A301").Select()
